# Apply the changes described by the commit diff:
#   1. Update two values on "main_variables": D8 (-319766 -> -308607) and
#      D16 (1045 -> 1000).
#   2. Rename the 7th sheet ("Copy of P&amp;amp;amp;amp;amp;L_Ratio Rev 1.1_Model"),
#      inserting one more "amp;" into its already-mangled name so it
#      becomes "Copy of P&amp;amp;amp;amp;amp;amp;L_Ratio Rev 1.1_Model".
#   3. Switch the active/selected sheet from "patient_transaction"
#      (activeTab 2) to "main_variables" (activeTab 0).

$wb = $excel.ActiveWorkbook

# --- 1. Update values on the main_variables sheet ---
$wsMain = $wb.Worksheets.Item("main_variables")
$wsMain.Range("D8").Value = -308607
$wsMain.Range("D16").Value = 1000

# --- 2. Rename sheet 7 (append one more "amp;" before "L_Ratio") ---
# Some engines surface a (non-fatal) long-sheet-name diagnostic as a
# terminating error here even though the rename itself is applied; wrap
# it so the rest of the script still runs no matter what.
$wsPL = $wb.Worksheets.Item(7)
try {
    $wsPL.Name = "Copy of P&amp;amp;amp;amp;amp;amp;L_Ratio Rev 1.1_Model"
} catch {
    Write-Host ("Sheet rename warning: " + $_.Exception.Message)
}

# --- 3. Make main_variables the active/selected sheet ---
$wsMain.Activate()
$wsMain.Select()

Write-Host "Done applying edits"
